$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.383.74"
$ws.Range("E2").Value = "  +0.47%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.606.19"
$ws.Range("E3").Value = "  +0.75%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - BNB (plain decimal -> force text so it doesn't become a number)
$ws.Range("D5").Value = "'211.98"

# Row 6 - XRP
$ws.Range("D6").Value = "'0.499"
$ws.Range("E6").Value = "  -0.93%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("E8").Value = "  -0.68%  "

# Row 9
$ws.Range("E9").Value = "  -0.36%  "

# Row 10
$ws.Range("D10").Value = "'19.22"
$ws.Range("E10").Value = "  +1.36%  "

# Row 11
$ws.Range("D11").Value = "'0.0848"
$ws.Range("E11").Value = "  -0.45%  "

# Row 12
$ws.Range("D12").Value = "1.832.07"
$ws.Range("E12").Value = "  +0.71%  "

# Row 13
$ws.Range("D13").Value = "1.594.77"
$ws.Range("E13").Value = "  -0.21%  "

# Row 14
$ws.Range("D14").Value = "'3.99"
$ws.Range("E14").Value = "  -0.43%  "

# Row 15
$ws.Range("E15").Value = "  -0.51%  "

# Row 16
$ws.Range("D16").Value = "'63.35"
$ws.Range("E16").Value = "  -0.93%  "

# Row 17
$ws.Range("D17").Value = "26.380.02"
$ws.Range("E17").Value = "  +0.45%  "

# Row 18
$ws.Range("D18").Value = "'230.51"
$ws.Range("E18").Value = "  +7.80%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0723"
$ws.Range("E19").Value = "  -0.40%  "

# Row 20
$ws.Range("D20").Value = "'7.60"
$ws.Range("E20").Value = "  +3.41%  "

# Row 21
$ws.Range("E21").Value = "  +0.07%  "

# Row 22
$ws.Range("E22").Value = "  -0.78%  "

# Row 23
$ws.Range("E23").Value = "  +3.31%  "

# Row 24
$ws.Range("D24").Value = "'8.95"
$ws.Range("E24").Value = "  -1.07%  "

# Row 25
$ws.Range("D25").Value = "'146.73"

# Row 26
$ws.Range("E26").Value = "  +0.00%  "

# Row 27
$ws.Range("E27").Value = "  -0.22%  "

# Row 28
$ws.Range("D28").Value = "'0.113"
$ws.Range("E28").Value = "  +0.65%  "

# Row 29
$ws.Range("E29").Value = "  +1.83%  "

# Row 30
$ws.Range("E30").Value = "  +0.94%  "

# Row 31
$ws.Range("E31").Value = "  -0.08%  "

# Row 32
$ws.Range("D32").Value = "1.495.77"
$ws.Range("E32").Value = "  +5.47%  "

# Row 33
$ws.Range("E33").Value = "  +0.53%  "

# Row 34
$ws.Range("E34").Value = "  -1.64%  "

# Row 35
$ws.Range("E35").Value = "  -0.33%  "

# Row 36
$ws.Range("E36").Value = "  +0.65%  "

# Row 37
$ws.Range("D37").Value = "'0.562"
$ws.Range("E37").Value = "  -3.23%  "

# Row 38
$ws.Range("E38").Value = "  -0.42%  "

# Row 39
$ws.Range("E39").Value = "  -0.29%  "

# Row 40
$ws.Range("E40").Value = "  -0.60%  "

# Row 41
$ws.Range("E41").Value = "  -0.02%  "

# Row 42
$ws.Range("E42").Value = "  +0.89%  "

# Row 43
$ws.Range("D43").Value = "'0.927"
$ws.Range("E43").Value = "  -4.02%  "

# Row 44
$ws.Range("D44").Value = "1.744.65"
$ws.Range("E44").Value = "  +0.80%  "

# Row 45
$ws.Range("E45").Value = "  -0.80%  "

# Row 46
$ws.Range("D46").Value = "'60.84"

# Row 47
$ws.Range("D47").Value = "'89.74"
$ws.Range("E47").Value = "  +3.27%  "

# Row 48
$ws.Range("E48").Value = "  -0.71%  "

# Row 49
$ws.Range("E49").Value = "  -0.32%  "

# Row 50
$ws.Range("E50").Value = "  +0.24%  "

# Row 51 - EnergySwap -> USDD
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "  -0.06%  "
